$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 4; rows 4-30 shift down to 5-31.
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the new record.
$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C4").Value = "Arica y Parinacota"
$ws.Range("D4").Value = 44545
$ws.Range("E4").Value = 15
$ws.Range("F4").Value = "Fruta"
$ws.Range("G4").Value = 100103
$ws.Range("H4").Value = "Frutos de hueso (carozo)"
$ws.Range("I4").Value = 100103004
$ws.Range("J4").Value = "Durazno"
$ws.Range("K4").Value = "Early Majestic"
$ws.Range("L4").Value = "Segunda"
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 19000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 19500
$ws.Range("Q4").Value = "$/bandeja 18 kilos granel"
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1083
$ws.Range("T4").Value = 18
